$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data

$ws.Range("D2").Value = '64.110.93'
$ws.Range("E2").Value = '  -0.34%  '

$ws.Range("D3").Value = '3.477.31'
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.24'
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.44'
$ws.Range("E6").Value = '  -2.25%  '

$ws.Range("E8").Value = '  -0.98%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.61'
$ws.Range("E9").Value = '  +4.95%  '

$ws.Range("E10").Value = '  -1.92%  '

$ws.Range("E11").Value = '  -0.61%  '

$ws.Range("D12").Value = '4.070.99'
$ws.Range("E12").Value = '  -0.59%  '

$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("E14").Value = '  -2.95%  '

$ws.Range("D15").Value = '3.478.20'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = '64.127.67'
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("E17").Value = '  -6.58%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.98'
$ws.Range("E18").Value = '  +0.52%  '

$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.41'
$ws.Range("E20").Value = '  -2.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '384.46'
$ws.Range("E21").Value = '  -2.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.570'
$ws.Range("E22").Value = '  -0.49%  '

$ws.Range("D23").Value = '3.617.72'
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.79'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("E26").Value = '  -0.54%  '

$ws.Range("E27").Value = '  -2.72%  '

$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  -3.86%  '

$ws.Range("E31").Value = '  -4.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.93'
$ws.Range("E32").Value = '  -4.79%  '

$ws.Range("D33").Value = '3.507.07'
$ws.Range("E33").Value = '  -0.42%  '

$ws.Range("E34").Value = '  +0.74%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.94'
$ws.Range("E36").Value = '  -2.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.75'
$ws.Range("E38").Value = '  -2.60%  '

$ws.Range("E39").Value = '  -4.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.64'
$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("E41").Value = '  -1.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.797'
$ws.Range("E42").Value = '  -1.25%  '

$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("E44").Value = '  -1.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.29'
$ws.Range("E45").Value = '  -3.13%  '

$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.30'
$ws.Range("E47").Value = '  -7.57%  '

$ws.Range("E48").Value = '  -4.05%  '

$ws.Range("E49").Value = '  -1.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.900'
$ws.Range("E50").Value = '  +0.43%  '

$ws.Range("D51").Value = '2.351.30'
$ws.Range("E51").Value = '  -4.76%  '
